$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.661.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "'2.801.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'355.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "'109.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "'0.557"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.10%  "
$ws.Range("D10").Value = "'39.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "'0.0836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "'20.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "'3.240.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "'2.816.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "'0.947"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'51.619.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +3.71%  "
$ws.Range("D20").Value = "'3.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("D21").Value = "'13.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.40%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'70.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'267.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "'2.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'26.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").Value = "'10.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "'37.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.19%  "
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'51.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "'5.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.76%  "
$ws.Range("D35").Value = "'0.0445"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.46%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'18.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("D43").Value = "'119.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "'21.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "'2.125.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").Value = "'2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.19%  "
$ws.Range("D48").Value = "'3.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "'0.905"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.62%  "
$ws.Range("D50").Value = "'5.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.18%  "
$ws.Range("E51").Value = "  +6.97%  "
